# Apply cryptos list update (GitHub Actions scrape refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-44, 48-51: refreshed Price (D) / Volume(1h) (E) values.
# Leading apostrophe forces text, matching the original inlineStr cell
# type (prevents Excel's automatic number/percent coercion).
$ws.Range("D2").Value = "'63.408.80"
$ws.Range("E2").Value = "'  -2.15%  "
$ws.Range("D3").Value = "'3.123.91"
$ws.Range("E3").Value = "'  -3.48%  "
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'560.49"
$ws.Range("E5").Value = "'  -3.17%  "
$ws.Range("D6").Value = "'161.30"
$ws.Range("E6").Value = "'  -8.12%  "
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("D8").Value = "'0.580"
$ws.Range("E8").Value = "'  -7.98%  "
$ws.Range("D9").Value = "'6.67"
$ws.Range("E9").Value = "'  -1.31%  "
$ws.Range("D10").Value = "'0.114"
$ws.Range("E10").Value = "'  -6.26%  "
$ws.Range("D11").Value = "'0.377"
$ws.Range("E11").Value = "'  -3.88%  "
$ws.Range("D12").Value = "'3.664.39"
$ws.Range("E12").Value = "'  -3.79%  "
$ws.Range("E13").Value = "'  -0.74%  "
$ws.Range("D14").Value = "'63.384.25"
$ws.Range("E14").Value = "'  -2.48%  "
$ws.Range("D15").Value = "'24.74"
$ws.Range("E15").Value = "'  -3.70%  "
$ws.Range("D16").Value = "'3.098.19"
$ws.Range("E16").Value = "'  -5.22%  "
$ws.Range("D17").Value = "'0.0000152"
$ws.Range("E17").Value = "'  -4.13%  "
$ws.Range("D18").Value = "'402.57"
$ws.Range("E18").Value = "'  -3.16%  "
$ws.Range("D19").Value = "'5.18"
$ws.Range("E19").Value = "'  -3.64%  "
$ws.Range("D20").Value = "'12.35"
$ws.Range("E20").Value = "'  -3.88%  "
$ws.Range("D21").Value = "'7.00"
$ws.Range("E21").Value = "'  -2.46%  "
$ws.Range("E22").Value = "'  +0.02%  "
$ws.Range("D23").Value = "'67.05"
$ws.Range("E23").Value = "'  -4.69%  "
$ws.Range("D24").Value = "'0.197"
$ws.Range("E24").Value = "'  -2.54%  "
$ws.Range("D25").Value = "'0.477"
$ws.Range("E25").Value = "'  -3.87%  "
$ws.Range("E26").Value = "'  -8.42%  "
$ws.Range("D27").Value = "'8.68"
$ws.Range("E27").Value = "'  -5.63%  "
$ws.Range("E28").Value = "'  -0.12%  "
$ws.Range("E29").Value = "'  -0.07%  "
$ws.Range("D30").Value = "'1.78"
$ws.Range("E30").Value = "'  -5.11%  "
$ws.Range("D31").Value = "'20.83"
$ws.Range("E31").Value = "'  -4.61%  "
$ws.Range("D32").Value = "'6.18"
$ws.Range("E32").Value = "'  -3.92%  "
$ws.Range("D33").Value = "'4.71"
$ws.Range("E33").Value = "'  -6.77%  "
$ws.Range("E34").Value = "'  -4.54%  "
$ws.Range("D35").Value = "'152.81"
$ws.Range("E35").Value = "'  -2.71%  "
$ws.Range("E36").Value = "'  -6.69%  "
$ws.Range("D37").Value = "'2.738.86"
$ws.Range("E37").Value = "'  -2.85%  "
$ws.Range("E38").Value = "'  -5.86%  "
$ws.Range("D39").Value = "'23.22"
$ws.Range("E39").Value = "'  -9.19%  "
$ws.Range("D40").Value = "'4.01"
$ws.Range("E40").Value = "'  -4.99%  "
$ws.Range("E41").Value = "'  -4.67%  "
$ws.Range("D42").Value = "'0.0613"
$ws.Range("E42").Value = "'  -2.16%  "
$ws.Range("E43").Value = "'  -6.97%  "
$ws.Range("D44").Value = "'0.0255"
$ws.Range("E44").Value = "'  -3.06%  "
$ws.Range("D48").Value = "'0.0968"
$ws.Range("E48").Value = "'  -4.10%  "
$ws.Range("E49").Value = "'  +0.78%  "
$ws.Range("D50").Value = "'1.87"
$ws.Range("E50").Value = "'  -13.99%  "
$ws.Range("D51").Value = "'5.65"
$ws.Range("E51").Value = "'  -3.04%  "

# Rows 45-47: coin ranking re-sorted (FirstDigitalUSD moved up to #45,
# Bittensor to #46, InjectiveProtocol to #47) with refreshed data.
$ws.Range("B45").Value = "'FirstDigitalUSD"
$ws.Range("C45").Value = "'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "'  +0.05%  "
$ws.Range("B46").Value = "'Bittensor"
$ws.Range("C46").Value = "'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").Value = "'279.33"
$ws.Range("E46").Value = "'  -8.56%  "
$ws.Range("B47").Value = "'InjectiveProtocol"
$ws.Range("C47").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'20.53"
$ws.Range("E47").Value = "'  -7.94%  "
